$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

# Insert a new row at position 19, pushing existing rows (19-40) down to (20-41)
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row with the new key/value pair
$ws.Cells.Item(19, 1).Value = "new_house"
$ws.Cells.Item(19, 2).Value = "New house available! Deploy it to continue."

# Match formatting used by the rest of the data rows (column B uses style index 2 = wrap text)
$ws.Cells.Item(20, 2).Copy()
$ws.Cells.Item(19, 2).PasteSpecial(-4122)

# Update selection to match the target state
$ws.Range("B19").Select()
